# Fix the truncated Pos value for the 10002 row (D4): "...448,172" -> "...448,1728"
# and nudge the active-cell selection on Sheet1 from D8 to D7, matching the
# author's commit ("add alot asset and add card bag and achieve random
# generation for the normal level").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D4").Value = "2040,852,1728=2040,448,1728"

$ws.Range("D7").Select()
